# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.357.17"
$ws.Range("E2").Value = "  +0.14%  "

# Row 3
$ws.Range("D3").Value = "1.867.12"
$ws.Range("E3").Value = "  +3.25%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9979"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.45%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.54"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.65%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9982"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.37%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5000"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.88%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3981"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.16%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1012"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +28.87%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.122"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.85%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.34"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.75%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.480"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.24%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.92"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.38%  "

# Row 14
$ws.Range("D14").Value = "1.848.34"
$ws.Range("E14").Value = "  +2.02%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.397"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.23%  "

# Row 16
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9971"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.51%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.63"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.15%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06653"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.48%  "

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.39"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.50%  "

# Row 21
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9976"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.42%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.063"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.77%  "

# Row 23
$ws.Range("D23").Value = "28.452.62"
$ws.Range("E23").Value = "  +0.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.36"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.11%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.250"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.94%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.17"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.15%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.485"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.76%  "

# Row 28
$ws.Range("D28").Value = "2.059.25"
$ws.Range("E28").Value = "  +1.89%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.66"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.10%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.38"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.39%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1057"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.12%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.054"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.93%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.653"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.50%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.604"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.82%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06812"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.34%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.144"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.16%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02388"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.48%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2167"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.92%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.033"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.34%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.53"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.49%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6290"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.71%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.181"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.99%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9982"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.36%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.34"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.82%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5995"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.15%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.695"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.45%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.280"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.67%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.92"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.35%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.976"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.85%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.192"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.25%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06849"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.18%  "
